{"js": "const body = context.document.body;\n\n// Ordered list of (old formula text, new formula text) pairs, taken from the\n// OOXML diff. Every old value is a unique substring in the document, so an\n// exact, case-sensitive search reliably finds the single matching run.\nconst replacements = [\n  [\"27+44=\", \"35+17=\"],\n  [\"37+2=\", \"95-5=\"],\n  [\"64-3=\", \"48-26=\"],\n  [\"5+11=\", \"89-75=\"],\n  [\"45-35=\", \"16-1=\"],\n  [\"46+51=\", \"75-41=\"],\n  [\"43+20=\", \"15+56=\"],\n  [\"91-42=\", \"8+73=\"],\n  [\"52+33=\", \"51-23=\"],\n  [\"59+35=\", \"57-22=\"],\n  [\"25+13=\", \"7+71=\"],\n  [\"69+20=\", \"37-29=\"],\n  [\"18-15=\", \"10+57=\"],\n  [\"21-19=\", \"50-38=\"],\n  [\"74-12=\", \"7+83=\"],\n  [\"76-31=\", \"32+38=\"],\n  [\"13+85=\", \"62+16=\"],\n  [\"27+51=\", \"12-2=\"],\n  [\"5+24=\", \"57-50=\"],\n  [\"1+45=\", \"29+62=\"],\n  [\"24+1=\", \"84-13=\"],\n  [\"16+82=\", \"25+17=\"],\n  [\"6+42=\", \"21-0=\"],\n  [\"83-38=\", \"76-75=\"],\n  [\"26-14=\", \"17+61=\"],\n  [\"31-23=\", \"32-14=\"],\n  [\"42+43=\", \"32-24=\"],\n  [\"76-39=\", \"9-2=\"],\n  [\"25+64=\", \"57+9=\"],\n  [\"69-64=\", \"97-73=\"],\n  [\"61+20=\", \"30+51=\"],\n  [\"42-34=\", \"65-12=\"],\n  [\"14+44=\", \"99-98=\"],\n  [\"65-26=\", \"90-56=\"],\n  [\"60+9=\", \"19+23=\"],\n  [\"15+46=\", \"7+31=\"],\n  [\"4+65=\", \"1+64=\"],\n  [\"59+7=\", \"6+79=\"],\n  [\"85-80=\", \"78-73=\"],\n  [\"9+82=\", \"69+17=\"],\n  [\"54+26=\", \"47+50=\"],\n  [\"22+10=\", \"98-96=\"],\n  [\"77-54=\", \"6+31=\"],\n  [\"40-21=\", \"38+32=\"],\n  [\"1+90=\", \"82-35=\"],\n  [\"58-45=\", \"21+58=\"],\n  [\"77-25=\", \"50-10=\"],\n  [\"82-81=\", \"51+24=\"],\n  [\"73-48=\", \"54-19=\"],\n  [\"85-62=\", \"54+8=\"],\n  [\"41-26=\", \"81-31=\"],\n  [\"1+53=\", \"93-6=\"],\n  [\"6+84=\", \"12+44=\"],\n  [\"47+38=\", \"36+55=\"],\n  [\"20+29=\", \"75-34=\"],\n  [\"43+34=\", \"79-4=\"],\n  [\"96-79=\", \"6+52=\"],\n  [\"97-50=\", \"61+32=\"],\n  [\"14+27=\", \"58+8=\"],\n  [\"18+71=\", \"94-7=\"],\n  [\"28-3=\", \"64-53=\"],\n  [\"86-71=\", \"16-15=\"],\n  [\"9+20=\", \"54-44=\"],\n  [\"33+55=\", \"68-59=\"],\n  [\"31+14=\", \"34+53=\"],\n  [\"16+56=\", \"35+17=\"],\n  [\"28-11=\", \"79-7=\"],\n  [\"10+82=\", \"54-32=\"],\n  [\"89-65=\", \"38-16=\"],\n  [\"12+83=\", \"48-1=\"],\n  [\"26+46=\", \"90-33=\"],\n  [\"3+48=\", \"17-2=\"],\n  [\"9+89=\", \"17+34=\"],\n  [\"77-3=\", \"51-47=\"],\n  [\"89+10=\", \"82-40=\"],\n  [\"81-50=\", \"56-4=\"],\n  [\"65+10=\", \"99-38=\"],\n  [\"59-12=\", \"56+33=\"],\n  [\"20-14=\", \"55-47=\"],\n  [\"53-5=\", \"44+25=\"],\n  [\"69-0=\", \"48-4=\"],\n  [\"15+9=\", \"57+28=\"],\n  [\"52+40=\", \"63+32=\"],\n  [\"52-45=\", \"94-39=\"],\n  [\"89-57=\", \"51+27=\"],\n  [\"9-5=\", \"0+72=\"],\n  [\"37+20=\", \"85-67=\"],\n  [\"8+57=\", \"43+4=\"],\n  [\"52+22=\", \"98-73=\"],\n  [\"80-45=\", \"6+38=\"],\n  [\"53+18=\", \"2+86=\"],\n  [\"92-22=\", \"16+23=\"],\n  [\"93-37=\", \"85-71=\"],\n  [\"63-32=\", \"78-16=\"],\n  [\"15+60=\", \"64-12=\"],\n  [\"57-2=\", \"79-45=\"],\n  [\"39+56=\", \"70-58=\"],\n  [\"35-16=\", \"71-23=\"],\n  [\"56+16=\", \"13+84=\"],\n  [\"60-41=\", \"74-21=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  // eslint-disable-next-line no-await-in-loop\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly one match for \"${oldText}\" but found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n# Ordered list of (old formula text, new formula text) pairs taken from the OOXML\n# diff. Every old value is a unique substring of the document at the time it is\n# searched for (verified offline), so a bounded Find/Replace-all is safe and only\n# ever touches a single cell.\n$replacements = @(\n    ,@(\"27+44=\", \"35+17=\")\n    ,@(\"37+2=\", \"95-5=\")\n    ,@(\"64-3=\", \"48-26=\")\n    ,@(\"5+11=\", \"89-75=\")\n    ,@(\"45-35=\", \"16-1=\")\n    ,@(\"46+51=\", \"75-41=\")\n    ,@(\"43+20=\", \"15+56=\")\n    ,@(\"91-42=\", \"8+73=\")\n    ,@(\"52+33=\", \"51-23=\")\n    ,@(\"59+35=\", \"57-22=\")\n    ,@(\"25+13=\", \"7+71=\")\n    ,@(\"69+20=\", \"37-29=\")\n    ,@(\"18-15=\", \"10+57=\")\n    ,@(\"21-19=\", \"50-38=\")\n    ,@(\"74-12=\", \"7+83=\")\n    ,@(\"76-31=\", \"32+38=\")\n    ,@(\"13+85=\", \"62+16=\")\n    ,@(\"27+51=\", \"12-2=\")\n    ,@(\"5+24=\", \"57-50=\")\n    ,@(\"1+45=\", \"29+62=\")\n    ,@(\"24+1=\", \"84-13=\")\n    ,@(\"16+82=\", \"25+17=\")\n    ,@(\"6+42=\", \"21-0=\")\n    ,@(\"83-38=\", \"76-75=\")\n    ,@(\"26-14=\", \"17+61=\")\n    ,@(\"31-23=\", \"32-14=\")\n    ,@(\"42+43=\", \"32-24=\")\n    ,@(\"76-39=\", \"9-2=\")\n    ,@(\"25+64=\", \"57+9=\")\n    ,@(\"69-64=\", \"97-73=\")\n    ,@(\"61+20=\", \"30+51=\")\n    ,@(\"42-34=\", \"65-12=\")\n    ,@(\"14+44=\", \"99-98=\")\n    ,@(\"65-26=\", \"90-56=\")\n    ,@(\"60+9=\", \"19+23=\")\n    ,@(\"15+46=\", \"7+31=\")\n    ,@(\"4+65=\", \"1+64=\")\n    ,@(\"59+7=\", \"6+79=\")\n    ,@(\"85-80=\", \"78-73=\")\n    ,@(\"9+82=\", \"69+17=\")\n    ,@(\"54+26=\", \"47+50=\")\n    ,@(\"22+10=\", \"98-96=\")\n    ,@(\"77-54=\", \"6+31=\")\n    ,@(\"40-21=\", \"38+32=\")\n    ,@(\"1+90=\", \"82-35=\")\n    ,@(\"58-45=\", \"21+58=\")\n    ,@(\"77-25=\", \"50-10=\")\n    ,@(\"82-81=\", \"51+24=\")\n    ,@(\"73-48=\", \"54-19=\")\n    ,@(\"85-62=\", \"54+8=\")\n    ,@(\"41-26=\", \"81-31=\")\n    ,@(\"1+53=\", \"93-6=\")\n    ,@(\"6+84=\", \"12+44=\")\n    ,@(\"47+38=\", \"36+55=\")\n    ,@(\"20+29=\", \"75-34=\")\n    ,@(\"43+34=\", \"79-4=\")\n    ,@(\"96-79=\", \"6+52=\")\n    ,@(\"97-50=\", \"61+32=\")\n    ,@(\"14+27=\", \"58+8=\")\n    ,@(\"18+71=\", \"94-7=\")\n    ,@(\"28-3=\", \"64-53=\")\n    ,@(\"86-71=\", \"16-15=\")\n    ,@(\"9+20=\", \"54-44=\")\n    ,@(\"33+55=\", \"68-59=\")\n    ,@(\"31+14=\", \"34+53=\")\n    ,@(\"16+56=\", \"35+17=\")\n    ,@(\"28-11=\", \"79-7=\")\n    ,@(\"10+82=\", \"54-32=\")\n    ,@(\"89-65=\", \"38-16=\")\n    ,@(\"12+83=\", \"48-1=\")\n    ,@(\"26+46=\", \"90-33=\")\n    ,@(\"3+48=\", \"17-2=\")\n    ,@(\"9+89=\", \"17+34=\")\n    ,@(\"77-3=\", \"51-47=\")\n    ,@(\"89+10=\", \"82-40=\")\n    ,@(\"81-50=\", \"56-4=\")\n    ,@(\"65+10=\", \"99-38=\")\n    ,@(\"59-12=\", \"56+33=\")\n    ,@(\"20-14=\", \"55-47=\")\n    ,@(\"53-5=\", \"44+25=\")\n    ,@(\"69-0=\", \"48-4=\")\n    ,@(\"15+9=\", \"57+28=\")\n    ,@(\"52+40=\", \"63+32=\")\n    ,@(\"52-45=\", \"94-39=\")\n    ,@(\"89-57=\", \"51+27=\")\n    ,@(\"9-5=\", \"0+72=\")\n    ,@(\"37+20=\", \"85-67=\")\n    ,@(\"8+57=\", \"43+4=\")\n    ,@(\"52+22=\", \"98-73=\")\n    ,@(\"80-45=\", \"6+38=\")\n    ,@(\"53+18=\", \"2+86=\")\n    ,@(\"92-22=\", \"16+23=\")\n    ,@(\"93-37=\", \"85-71=\")\n    ,@(\"63-32=\", \"78-16=\")\n    ,@(\"15+60=\", \"64-12=\")\n    ,@(\"57-2=\", \"79-45=\")\n    ,@(\"39+56=\", \"70-58=\")\n    ,@(\"35-16=\", \"71-23=\")\n    ,@(\"56+16=\", \"13+84=\")\n    ,@(\"60-41=\", \"74-21=\")\n)\n\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = $wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($oldText, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceOne)\n    if (-not $found) {\n        throw \"Could not find expected text: $oldText\"\n    }\n}"}
